# Sync attendance_reports: swap order of "Recorded By" names in column G
# from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
    }
}
